# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns on the
# sheet with refreshed quote data, matching the commit's scraped snapshot.
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as literal text (e.g. "1.001", "0.07370") instead of
# auto-converting/normalizing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.207.52"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.853.90"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'235.35"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -2.78%  "
$ws.Range("D8").Value = "'0.2806"
$ws.Range("E8").Value = "  -4.06%  "
$ws.Range("D9").Value = "'0.06456"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("D10").Value = "1.855.46"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'0.07370"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'16.22"
$ws.Range("E12").Value = "  -4.05%  "
$ws.Range("D13").Value = "'5.087"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "'87.09"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "'0.6449"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "30.155.33"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'13.12"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "'0.000007560"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").Value = "'226.94"
$ws.Range("E20").Value = "  +18.34%  "
$ws.Range("D21").Value = "2.098.18"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'5.284"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'6.077"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'9.196"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("D26").Value = "'163.70"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "'18.46"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "'1.921"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "'1.439"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "'0.09178"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "'4.232"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "'3.954"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").Value = "'0.04964"
$ws.Range("E33").Value = "  -3.94%  "
$ws.Range("D34").Value = "'0.7300"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").Value = "'1.141"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "'0.01843"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").Value = "'2.593"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "'0.8988"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "'5.943"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").Value = "'105.87"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'0.4230"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").Value = "'7.358"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").Value = "'0.1310"
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("D47").Value = "'64.05"
$ws.Range("E47").Value = "  -6.45%  "
$ws.Range("D48").Value = "'1.495"
$ws.Range("E48").Value = "  +6.81%  "
$ws.Range("D49").Value = "'8.734"
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("D50").Value = "'33.77"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").Value = "'0.05652"
$ws.Range("E51").Value = "  -3.36%  "
